$wb = $excel.ActiveWorkbook

# --- Sheet "general" ---
$ws = $wb.Worksheets.Item("general")
$ws.Cells.Item(3, 2).Value = 441.4078650152845    # B3 objValue
$ws.Cells.Item(4, 2).Value = 0.0130000114440918   # B4 runtime
$ws.Cells.Item(6, 2).Value = 45.77786501528455    # B6 Z1
$ws.Cells.Item(10, 2).Value = 395.63              # B10 Z5

# --- Sheet "x" ---
$ws = $wb.Worksheets.Item("x")
$ws.Cells.Item(5, 2).Value = 11   # B5
$ws.Cells.Item(8, 2).Value = 10   # B8
$ws.Cells.Item(11, 2).Value = 13  # B11
$ws.Cells.Item(12, 2).Value = 12  # B12
$ws.Cells.Item(13, 2).Value = 7   # B13
$ws.Cells.Item(14, 2).Value = 1   # B14

# --- Sheet "U" ---
$ws = $wb.Worksheets.Item("U")
$ws.Cells.Item(3, 2).Value = 2   # B3

# --- Sheet "TBar" ---
$ws = $wb.Worksheets.Item("TBar")
$ws.Cells.Item(3, 2).Value = 22.61068956408072    # B3
$ws.Cells.Item(4, 2).Value = 10                   # B4
$ws.Cells.Item(9, 2).Value = 23.22876137241512    # B9
$ws.Cells.Item(11, 2).Value = 20                  # B11
$ws.Cells.Item(12, 2).Value = 20                  # B12
$ws.Cells.Item(13, 2).Value = 30                  # B13
$ws.Cells.Item(14, 2).Value = 25.21630137166173   # B14
$ws.Cells.Item(15, 2).Value = 28.21630585843227   # B15

# --- Sheet "Q" ---
$ws = $wb.Worksheets.Item("Q")
$ws.Cells.Item(7, 3).Value = 297.2549999999998    # C7
$ws.Cells.Item(8, 3).Value = 314.6649999999998    # C8
$ws.Cells.Item(9, 3).Value = 311.2049999999998    # C9
$ws.Cells.Item(10, 3).Value = 316.3799999999999   # C10
$ws.Cells.Item(11, 3).Value = 301.5249999999997   # C11
$ws.Cells.Item(12, 3).Value = 67.77500000000072   # C12
$ws.Cells.Item(13, 3).Value = 73.77000000000072   # C13
$ws.Cells.Item(14, 3).Value = 74.03500000000074   # C14
$ws.Cells.Item(15, 3).Value = 73.04500000000073   # C15
$ws.Cells.Item(16, 3).Value = 72.66500000000073   # C16
$ws.Cells.Item(25, 3).Value = 136.5699999999998   # C25
$ws.Cells.Item(37, 3).Value = 236.3350000000021   # C37
$ws.Cells.Item(38, 3).Value = 246.4550000000021   # C38
$ws.Cells.Item(39, 3).Value = 231.7250000000021   # C39
$ws.Cells.Item(40, 3).Value = 253.5450000000021   # C40
$ws.Cells.Item(41, 3).Value = 239.25              # C41
$ws.Cells.Item(47, 3).Value = 153.2600000000012   # C47
$ws.Cells.Item(48, 3).Value = 161.7350000000012   # C48
$ws.Cells.Item(49, 3).Value = 153.75              # C49
$ws.Cells.Item(50, 3).Value = 163.7750000000012   # C50
$ws.Cells.Item(51, 3).Value = 157.3950000000012   # C51
$ws.Cells.Item(52, 3).Value = 187.309999999999    # C52
$ws.Cells.Item(53, 3).Value = 197.074999999999    # C53
$ws.Cells.Item(54, 3).Value = 197.934999999999    # C54
$ws.Cells.Item(55, 3).Value = 197.2249999999991   # C55
$ws.Cells.Item(56, 3).Value = 185.719999999999    # C56
$ws.Cells.Item(57, 3).Value = 297.2549999999998   # C57
$ws.Cells.Item(58, 3).Value = 314.6649999999998   # C58
$ws.Cells.Item(59, 3).Value = 311.2049999999998   # C59
$ws.Cells.Item(60, 3).Value = 316.3799999999999   # C60
$ws.Cells.Item(61, 3).Value = 301.5249999999997   # C61
$ws.Cells.Item(62, 3).Value = 212.0549999999987   # C62
$ws.Cells.Item(63, 3).Value = 215.8299999999987   # C63
$ws.Cells.Item(64, 3).Value = 177.0399999999987   # C64
$ws.Cells.Item(65, 3).Value = 198                 # C65
$ws.Cells.Item(66, 3).Value = 184.7               # C66
$ws.Cells.Item(67, 3).Value = 236.3350000000021   # C67
$ws.Cells.Item(68, 3).Value = 246.4550000000021   # C68
$ws.Cells.Item(69, 3).Value = 231.7250000000021   # C69
$ws.Cells.Item(70, 3).Value = 253.5450000000021   # C70
$ws.Cells.Item(71, 3).Value = 239.25              # C71

# --- Sheet "R" ---
$ws = $wb.Worksheets.Item("R")
$ws.Cells.Item(7, 3).Value = 13.7   # C7
$ws.Cells.Item(8, 3).Value = 6.91   # C8
$ws.Cells.Item(9, 3).Value = 10.68  # C9
$ws.Cells.Item(10, 3).Value = 7.39  # C10
$ws.Cells.Item(11, 3).Value = 14.68 # C11
